# map_24_INIT_Positions.xlsx edit
# - Row 1 (C1:BJ1) index strip: was descending 59..0, now ascending 0..59
# - BK1 / A2 axis labels "Index 2" / "Index 1 " swapped -> "Index 1" / "Index 2"
# - BO3:BO12 "start point" sample column re-rolled with new values
#   (BQ3:BQ12 hold "=BOx*0.05" formulas and recalc automatically)
# - Selection moved to BQ14 (best-effort; matches author's final cursor spot)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: column index strip, now ascending 0..59 across C1:BJ1 ---
for ($col = 3; $col -le 62; $col++) {
    $ws.Cells.Item(1, $col).Value = $col - 3
}

# --- Axis labels: swap the two "Index" headers ---
$ws.Range("BK1").Value = "Index 1"
$ws.Range("A2").Value = "Index 2"

# --- New training run sample data (column BO), rows 3-12 ---
$ws.Range("BO3").Value = 48
$ws.Range("BO4").Value = 43
$ws.Range("BO5").Value = 23
$ws.Range("BO6").Value = 9
$ws.Range("BO7").Value = 24
$ws.Range("BO8").Value = 21
$ws.Range("BO9").Value = 6
$ws.Range("BO10").Value = 26
$ws.Range("BO11").Value = 41
$ws.Range("BO12").Value = 48

# --- Cursor / selection ---
$null = $ws.Range("BQ14").Select()
